$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E on this sheet hold free-form text (prices with
# "."-group separators, percentages with padding spaces). For any
# Price-column (D) value that looks like a plain decimal (e.g.
# "18.22"), Excel would otherwise silently coerce the assignment
# into a number. Force the cell to Text format first, write the
# string, then restore the original (default) cell style so no
# stray style index/attribute is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.926.38'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.623.95'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.503'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.90%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('E8').Value = '  -2.34%  '
$ws.Range('E9').Value = '  -3.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.22'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.82%  '
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.848.92'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.08%  '
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.619.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('E15').Value = '  -3.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.916.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.51%  '
$ws.Range('E18').Value = '  -3.76%  '
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.48%  '
$ws.Range('E21').Value = '  -2.93%  '
$ws.Range('E22').Value = '  -3.68%  '
$ws.Range('E23').Value = '  -2.34%  '
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.75'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.47%  '
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('E31').Value = '  -2.73%  '
$ws.Range('E32').Value = '  -3.99%  '
$ws.Range('E33').Value = '  -5.40%  '
$ws.Range('E34').Value = '  -2.93%  '
$ws.Range('E35').Value = '  -2.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.117.05'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.10%  '
$ws.Range('E37').Value = '  -6.53%  '
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.520'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.48%  '
$ws.Range('E40').Value = '  -2.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('E42').Value = '  -3.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.760.05'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('E44').Value = '  -5.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₆0115'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0528'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.49%  '
$ws.Range('E48').Value = '  -1.12%  '
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.37%  '
